# Updated cryptos list on Mon Feb 26 20:52:13 UTC 2024 with GitHub Actions
# Refreshes the Price (column D) and Volume(1h) (column E) columns with the
# latest scraped values. Numeric-looking price strings are written with a
# leading apostrophe so Excel keeps them as text (matching the source data,
# which stores prices as plain strings, not numbers) instead of silently
# reformatting/rounding them as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.611.68"
$ws.Range("E2").Value = "  +5.49%  "
$ws.Range("D3").Value = "3.188.69"
$ws.Range("E3").Value = "  +2.76%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'402.31"
$ws.Range("E5").Value = "  +3.74%  "
$ws.Range("D6").Value = "'109.17"
$ws.Range("E6").Value = "  +5.34%  "
$ws.Range("D7").Value = "'0.551"
$ws.Range("E7").Value = "  +1.21%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("E9").Value = "  +5.37%  "
$ws.Range("D10").Value = "'39.00"
$ws.Range("E10").Value = "  +5.09%  "
$ws.Range("E11").Value = "  +1.73%  "
$ws.Range("E12").Value = "  +2.20%  "
$ws.Range("D13").Value = "3.675.37"
$ws.Range("E13").Value = "  +2.42%  "
$ws.Range("D14").Value = "'19.25"
$ws.Range("E14").Value = "  +2.75%  "
$ws.Range("E15").Value = "  +2.82%  "
$ws.Range("E16").Value = "  +8.78%  "
$ws.Range("D17").Value = "3.173.18"
$ws.Range("E17").Value = "  +2.20%  "
$ws.Range("E18").Value = "  -1.92%  "
$ws.Range("D19").Value = "54.447.10"
$ws.Range("E19").Value = "  +4.86%  "
$ws.Range("E20").Value = "  +4.05%  "
$ws.Range("D21").Value = "'12.87"
$ws.Range("E21").Value = "  +2.72%  "
$ws.Range("D22").Value = "0.0₃0997"
$ws.Range("E22").Value = "  +2.68%  "
$ws.Range("D23").Value = "'71.73"
$ws.Range("E23").Value = "  +2.18%  "
$ws.Range("D24").Value = "'274.56"
$ws.Range("E24").Value = "  +2.07%  "
$ws.Range("D25").Value = "'3.28"
$ws.Range("E25").Value = "  +4.57%  "
$ws.Range("D26").Value = "'8.00"
$ws.Range("E26").Value = "  -2.38%  "
$ws.Range("D27").Value = "'27.79"
$ws.Range("E27").Value = "  +2.42%  "
$ws.Range("D28").Value = "'7.40"
$ws.Range("E28").Value = "  +2.17%  "
$ws.Range("E29").Value = "  -1.32%  "
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("D31").Value = "'0.113"
$ws.Range("E31").Value = "  +3.88%  "
$ws.Range("D32").Value = "'11.10"
$ws.Range("E32").Value = "  +7.03%  "
$ws.Range("E33").Value = "  +10.62%  "
$ws.Range("D34").Value = "'36.78"
$ws.Range("E34").Value = "  +3.32%  "
$ws.Range("E35").Value = "  +0.47%  "
$ws.Range("D36").Value = "'50.75"
$ws.Range("E36").Value = "  +0.80%  "
$ws.Range("D37").Value = "'3.67"
$ws.Range("E37").Value = "  +8.54%  "
$ws.Range("D38").Value = "'0.999"
$ws.Range("E38").Value = "  -0.14%  "
$ws.Range("E39").Value = "  +9.56%  "
$ws.Range("D40").Value = "'4.11"
$ws.Range("E40").Value = "  +10.82%  "
$ws.Range("E41").Value = "  +2.61%  "
$ws.Range("E42").Value = "  -1.10%  "
$ws.Range("D43").Value = "'17.32"
$ws.Range("E43").Value = "  +1.55%  "
$ws.Range("D44").Value = "'130.01"
$ws.Range("E44").Value = "  +1.96%  "
$ws.Range("E45").Value = "  +1.31%  "
$ws.Range("D46").Value = "'22.32"
$ws.Range("E46").Value = "  +0.46%  "
$ws.Range("D47").Value = "'2.51"
$ws.Range("E47").Value = "  +1.88%  "
$ws.Range("E48").Value = "  -0.83%  "
$ws.Range("D49").Value = "2.089.76"
$ws.Range("E49").Value = "  +1.75%  "
$ws.Range("D50").Value = "'0.0344"
$ws.Range("E50").Value = "  +7.72%  "
$ws.Range("D51").Value = "'0.0507"
$ws.Range("E51").Value = "  +11.05%  "
